# Apply cryptos list price/volume update (GitHub Actions refresh).
# Note: some "Price" (column D) values are plain decimals (e.g. "245.32")
# which Excel would otherwise auto-coerce to a Number when assigned via
# .Value. The source cells are text, so those assignments are prefixed
# with a literal leading apostrophe (quote-prefix) to force text storage,
# matching the original inline-string cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.464.00'
$ws.Range("E2").Value = '  +1.08%  '

$ws.Range("D3").Value = '1.903.89'
$ws.Range("E3").Value = '  +2.47%  '

$ws.Range("E4").Value = '  +0.41%  '

$ws.Range("D5").Value = '''245.32'
$ws.Range("E5").Value = '  +4.06%  '

$ws.Range("D6").Value = '''0.633'
$ws.Range("E6").Value = '  +1.48%  '

$ws.Range("E7").Value = '  +0.35%  '

$ws.Range("D8").Value = '''41.97'
$ws.Range("E8").Value = '  -2.24%  '

$ws.Range("E9").Value = '  +2.79%  '

$ws.Range("E10").Value = '  +1.03%  '

$ws.Range("E11").Value = '  +0.91%  '

$ws.Range("D12").Value = '2.180.19'
$ws.Range("E12").Value = '  +2.54%  '

$ws.Range("D13").Value = '''12.29'
$ws.Range("E13").Value = '  +7.45%  '

$ws.Range("E14").Value = '  +1.68%  '

$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '''4.85'
$ws.Range("E15").Value = '  +3.29%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.857.03'
$ws.Range("E16").Value = '  -0.29%  '

$ws.Range("D17").Value = '35.509.83'
$ws.Range("E17").Value = '  +1.35%  '

$ws.Range("D18").Value = '''71.71'
$ws.Range("E18").Value = '  +1.82%  '

$ws.Range("D20").Value = '''243.15'
$ws.Range("E20").Value = '  +0.64%  '

$ws.Range("D21").Value = '''12.59'
$ws.Range("E21").Value = '  +3.06%  '

$ws.Range("E22").Value = '  +1.72%  '

$ws.Range("E24").Value = '  +0.33%  '

$ws.Range("D25").Value = '''172.61'
$ws.Range("E25").Value = '  +0.65%  '

$ws.Range("D26").Value = '''2.19'
$ws.Range("E26").Value = '  +19.59%  '

$ws.Range("D27").Value = '''8.54'
$ws.Range("E27").Value = '  +7.74%  '

$ws.Range("D28").Value = '''17.94'
$ws.Range("E28").Value = '  +1.41%  '

$ws.Range("E29").Value = '  +0.39%  '

$ws.Range("D30").Value = '''0.974'
$ws.Range("E30").Value = '  +24.30%  '

$ws.Range("E31").Value = '  +2.47%  '

$ws.Range("D32").Value = '''4.10'
$ws.Range("E32").Value = '  +2.26%  '

$ws.Range("E33").Value = '  +0.49%  '

$ws.Range("E34").Value = '  +3.72%  '

$ws.Range("E35").Value = '  +6.76%  '

$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").Value = '''2.03'
$ws.Range("E36").Value = '  -0.24%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '''1.35'
$ws.Range("E37").Value = '  +8.95%  '

$ws.Range("E38").Value = '  +1.93%  '

$ws.Range("E39").Value = '  +1.50%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '''0.0631'
$ws.Range("E40").Value = '  +16.31%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '''90.74'
$ws.Range("E41").Value = '  -1.18%  '

$ws.Range("D42").Value = '''15.69'
$ws.Range("E42").Value = '  +4.17%  '

$ws.Range("E43").Value = '  +45.10%  '

$ws.Range("D44").Value = '1.348.77'
$ws.Range("E44").Value = '  -0.21%  '

$ws.Range("E45").Value = '  +1.98%  '

$ws.Range("D46").Value = '''12.78'
$ws.Range("E46").Value = '  -0.73%  '

$ws.Range("E47").Value = '  +0.44%  '

$ws.Range("D48").Value = '''2.75'
$ws.Range("E48").Value = '  -0.67%  '

$ws.Range("E49").Value = '  +3.12%  '

$ws.Range("D50").Value = '2.091.07'
$ws.Range("E50").Value = '  +2.69%  '

$ws.Range("D51").Value = '''0.0692'
$ws.Range("E51").Value = '  +1.40%  '
